# The document contains two structured document tags (SDT elements, aka
# "content controls" / fields in the Word object model):
#   - a block-level one aliased "BlockField" with tag "BlockField"
#   - an inline one aliased "InlineField" with tag "InlineField"
#
# Re-tag them using the ContentControls collection - the special
# collection that exposes a document's SDT/field elements - instead of
# touching the visible placeholder text, since the tag is metadata, not
# document text reachable via Find/Replace.

$d = $word.ActiveDocument

$contentControls = $d.ContentControls

for ($i = 1; $i -le $contentControls.Count; $i++) {
    $cc = $contentControls.Item($i)

    if ($cc.Tag -eq "BlockField") {
        $cc.Tag = "FirstTag"
    }
    elseif ($cc.Tag -eq "InlineField") {
        $cc.Tag = "SecondTag"
    }
}
